$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: nuclear concentration correction (Pu239) ---
$ws.Range("C5").Value = 0.0021052
$ws.Range("D5").Value = 265.76
$ws.Range("E5").Value = 2.86

# --- Row 6: nuclear concentration correction (Pu241) ---
$ws.Range("C6").Value = 0.00068223
$ws.Range("D6").Value = 237.79

# --- Row 14: kinf / L^2 inputs (1 group) ---
$ws.Range("B14").Value = 0.4
$ws.Range("C14").Value = 0.02266236

# apply scientific number format to E14 (alpha^2), matching the style
# already used for similar small scientific-notation cells (C5/C6/C14/C18)
$ws.Range("E14").NumberFormat = "0.00E+00"

# --- Row 18: kinf / L^2 inputs (2 group) / start of tsunami temp rise calc ---
$ws.Range("B18").Value = 0.3792084
$ws.Range("C18").Value = 0.112136

# --- Update active selection to B11 ---
$ws.Range("B11").Select() | Out-Null
